$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 708, shifting existing rows 708:805 down to 709:806
$ws.Rows.Item(708).Insert()

# Populate the newly inserted row 708 with its data.
$ws.Cells.Item(708, 1).Value = 6
$ws.Cells.Item(708, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(708, 3).Value = "Metropolitana"
$ws.Cells.Item(708, 4).Value = 45154
$ws.Cells.Item(708, 5).Value = 13
$ws.Cells.Item(708, 6).Value = 100112012
$ws.Cells.Item(708, 7).Value = "Espinaca"
$ws.Cells.Item(708, 8).Value = "Sin especificar"
$ws.Cells.Item(708, 9).Value = "Primera"
$ws.Cells.Item(708, 10).Value = 550
$ws.Cells.Item(708, 11).Value = 5000
$ws.Cells.Item(708, 12).Value = 6000
$ws.Cells.Item(708, 13).Value = 5455
$ws.Cells.Item(708, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(708, 15).Value = "Región Metropolitana"
$ws.Cells.Item(708, 16).Value = 546
$ws.Cells.Item(708, 17).Value = 10
$ws.Cells.Item(708, 18).Value = "Hortaliza"

# Match the date style used by the rest of column D.
$ws.Cells.Item(708, 4).NumberFormat = $ws.Cells.Item(709, 4).NumberFormat
